$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-formatted style from the last existing data row (A59)
# down onto the two new rows so the new date cells pick up the same
# numFmtId-14 "m/d/yyyy"-style short-date format (style index 1) instead
# of Excel's General format.
$ws.Cells.Item(59, 1).Copy()
$ws.Cells.Item(60, 1).PasteSpecial(-4122)
$ws.Cells.Item(61, 1).PasteSpecial(-4122)

# New row 60: 1/16/2026 (serial 46038), Error Count 2
$ws.Cells.Item(60, 1).Value = 46038
$ws.Cells.Item(60, 2).Value = 2

# New row 61: 1/15/2026 (serial 46037), Error Count 9
$ws.Cells.Item(61, 1).Value = 46037
$ws.Cells.Item(61, 2).Value = 9

# Match the author's final on-screen selection covering both new rows.
$ws.Range("A60:B61").Select()
